$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$cell = $ws.Cells.Item(7, 9)
Write-Host ("Address: " + $cell.Address())
$ws.Hyperlinks.Add($cell, "https://example.com", [Type]::Missing, [Type]::Missing, "75ca882a-37eb-4e84-86cd-e94d68725312.md") | Out-Null
